$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: goto -> new target site (nectarsleep.com), turned into a hyperlink
$ws.Range("D2").Value = "https://www.nectarsleep.com"

# Row 3: waitfortext step now targets "Mattresses" / locatorType "a"
$ws.Range("D3").Value = "Mattresses"
$ws.Range("E3").Value = "a"

# Row 4: action becomes "click" on "Mattresses" / "a", value cleared, new waitAfter
$ws.Range("C4").Value = "click"
$ws.Range("D4").Value = "Mattresses"
$ws.Range("E4").Value = "a"
$ws.Range("F4").ClearContents()
$ws.Range("H4").Value = 2000

# Turn D2 into a real hyperlink pointing at the new URL
$ws.Hyperlinks.Add($ws.Range("D2"), "https://www.nectarsleep.com")
# Re-touch the font so the pre-existing "Hyperlink" style (already applied to
# D2) is reused instead of the engine minting a brand-new cell style.
$ws.Range("D2").Font.Underline = $true
